# Refresh crypto price/volume snapshot (GitHub Actions update).
# Values are written as literal text so formatting such as trailing
# zeros, thousands separators and exact decimal digits is preserved
# exactly as scraped, instead of being reinterpreted as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.959.20'
$ws.Range('E2').Value = '  +2.00%  '
$ws.Range('D3').Value = '1.846.61'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('D4').Value = '''1.009'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''310.11'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = '''0.4668'
$ws.Range('E7').Value = '  +3.21%  '
$ws.Range('D8').Value = '''0.3628'
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('D9').Value = '''0.07172'
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('D10').Value = '''0.9264'
$ws.Range('E10').Value = '  +4.15%  '
$ws.Range('D11').Value = '''19.60'
$ws.Range('E11').Value = '  +0.95%  '
$ws.Range('D12').Value = '''0.07678'
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('D13').Value = '1.855.05'
$ws.Range('E13').Value = '  +2.09%  '
$ws.Range('D14').Value = '''5.295'
$ws.Range('E14').Value = '  +0.13%  '
$ws.Range('D15').Value = '''6.405'
$ws.Range('E15').Value = '  +1.42%  '
$ws.Range('D16').Value = '''88.35'
$ws.Range('E16').Value = '  +3.45%  '
$ws.Range('D17').Value = '''1.010'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '''0.000008595'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '26.989.71'
$ws.Range('E20').Value = '  +2.01%  '
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('D22').Value = '''5.033'
$ws.Range('E22').Value = '  +1.37%  '
$ws.Range('D23').Value = '''10.63'
$ws.Range('E23').Value = '  +1.12%  '
$ws.Range('D24').Value = '''1.937'
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('D25').Value = '''152.55'
$ws.Range('E25').Value = '  -0.28%  '
$ws.Range('D26').Value = '''18.14'
$ws.Range('E26').Value = '  +1.98%  '
$ws.Range('D27').Value = '''2.049'
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('D28').Value = '''113.93'
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('D29').Value = '''4.931'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').Value = '''0.08860'
$ws.Range('E30').Value = '  +1.83%  '
$ws.Range('D31').Value = '''3.182'
$ws.Range('E31').Value = '  +2.04%  '
$ws.Range('D32').Value = '''2.836'
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').Value = '''1.177'
$ws.Range('E33').Value = '  +6.49%  '
$ws.Range('D34').Value = '''0.7457'
$ws.Range('E34').Value = '  +3.10%  '
$ws.Range('D35').Value = '''4.474'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').Value = '''1.087'
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('E37').Value = '  +2.43%  '
$ws.Range('D38').Value = '''0.01940'
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('D39').Value = '''0.05167'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('D40').Value = '''0.5149'
$ws.Range('E40').Value = '  +1.35%  '
$ws.Range('D41').Value = '''6.894'
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('D42').Value = '''0.1511'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '''10.63'
$ws.Range('E43').Value = '  +6.89%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '''8.173'
$ws.Range('E44').Value = '  +1.91%  '
$ws.Range('D45').Value = '''0.4711'
$ws.Range('E45').Value = '  +0.86%  '
$ws.Range('D46').Value = '''1.007'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('D47').Value = '''100.34'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '''1.602'
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('D49').Value = '''0.06050'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('D50').Value = '''64.39'
$ws.Range('E50').Value = '  +1.09%  '
$ws.Range('D51').Value = '''36.11'
$ws.Range('E51').Value = '  -0.10%  '
